$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (ALC)
$ws.Range("H15").Value = 5088.578
$ws.Range("I15").Value = 5088.578
$ws.Range("K15").Value = 15265.734
$ws.Range("M15").Value = -15096.734

# Row 40 (ALC)
$ws.Range("H40").Value = 1005.1905
$ws.Range("I40").Value = 997.58826
$ws.Range("J40").Value = 1037.5
$ws.Range("K40").Value = 997.58826
$ws.Range("L40").Value = 1037.5
$ws.Range("M40").Value = -822.58826
$ws.Range("N40").Value = -1387.5

# Row 98 (ALC)
$ws.Range("H98").Value = 3309.5
$ws.Range("I98").Value = 3285.3125
$ws.Range("J98").Value = 3503
$ws.Range("K98").Value = 3285.3125
$ws.Range("L98").Value = 3503
$ws.Range("M98").Value = -1787.3125
$ws.Range("N98").Value = -6499

# Row 112 (ALC)
$ws.Range("H112").Value = 83334610
$ws.Range("J112").Value = 83334610
$ws.Range("L112").Value = 250003830
$ws.Range("N112").Value = -250006046

# Row 122 (ALC)
$ws.Range("H122").Value = 3309.5
$ws.Range("I122").Value = 3285.3125
$ws.Range("J122").Value = 3503
$ws.Range("K122").Value = 9855.9375
$ws.Range("L122").Value = 10509
$ws.Range("M122").Value = -7405.9375
$ws.Range("N122").Value = -15409

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1961197
$ws.Range("I2").Value = 415.33334
$ws.Range("J2").Value = 3268384.8
$ws.Range("K2").Value = 415.33334
$ws.Range("L2").Value = 3268384.8
$ws.Range("M2").Value = -302.33334
$ws.Range("N2").Value = -3268610.8

# Row 32 (ARM)
$ws.Range("H32").Value = 11548.551
$ws.Range("I32").Value = 6565.3657
$ws.Range("J32").Value = 37087.375
$ws.Range("K32").Value = 6565.3657
$ws.Range("L32").Value = 37087.375
$ws.Range("M32").Value = -6278.3657
$ws.Range("N32").Value = -37661.375

# Row 116 (ARM)
$ws.Range("H116").Value = 1961197
$ws.Range("I116").Value = 415.33334
$ws.Range("J116").Value = 3268384.8
$ws.Range("K116").Value = 415.33334
$ws.Range("L116").Value = 3268384.8
$ws.Range("M116").Value = 1878.66666
$ws.Range("N116").Value = -3272972.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1961197
$ws.Range("I3").Value = 415.33334
$ws.Range("J3").Value = 3268384.8
$ws.Range("K3").Value = 415.33334
$ws.Range("L3").Value = 3268384.8
$ws.Range("M3").Value = -301.33334
$ws.Range("N3").Value = -3268612.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 2223.6
$ws.Range("I31").Value = 1671.1428
$ws.Range("J31").Value = 3512.6667
$ws.Range("K31").Value = 1671.1428
$ws.Range("L31").Value = 3512.6667
$ws.Range("M31").Value = -1376.1428
$ws.Range("N31").Value = -4102.6667

# Row 34 (CRP)
$ws.Range("H34").Value = 2223.6
$ws.Range("I34").Value = 1671.1428
$ws.Range("J34").Value = 3512.6667
$ws.Range("K34").Value = 1671.1428
$ws.Range("L34").Value = 3512.6667
$ws.Range("M34").Value = -1469.1428
$ws.Range("N34").Value = -3916.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 31 (CUL)
$ws.Range("H31").Value = 2000
$ws.Range("J31").Value = 2000
$ws.Range("L31").Value = 6000
$ws.Range("N31").Value = -6576

# Row 41 (CUL)
$ws.Range("H41").Value = 2343.5557
$ws.Range("I41").Value = 392
$ws.Range("J41").Value = 2587.5
$ws.Range("K41").Value = 1176
$ws.Range("L41").Value = 7762.5
$ws.Range("M41").Value = -838
$ws.Range("N41").Value = -8438.5

# Row 58 (CUL)
$ws.Range("H58").Value = 2995
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9256

# Row 68 (CUL)
$ws.Range("H68").Value = 819.6517
$ws.Range("I68").Value = 606.28
$ws.Range("J68").Value = 1093.2051
$ws.Range("K68").Value = 1818.84
$ws.Range("L68").Value = 3279.615299999999
$ws.Range("M68").Value = -1007.84
$ws.Range("N68").Value = -4901.615299999999

# Row 71 (CUL)
$ws.Range("H71").Value = 819.6517
$ws.Range("I71").Value = 606.28
$ws.Range("J71").Value = 1093.2051
$ws.Range("K71").Value = 5456.52
$ws.Range("L71").Value = 9838.845899999998
$ws.Range("M71").Value = -1400.52
$ws.Range("N71").Value = -17950.8459

# Row 105 (CUL)
$ws.Range("H105").Value = 454001440
$ws.Range("J105").Value = 454001440
$ws.Range("L105").Value = 1362004320
$ws.Range("N105").Value = -1362009562

# Row 107 (CUL)
$ws.Range("H107").Value = 28269.549
$ws.Range("J107").Value = 31511.484
$ws.Range("L107").Value = 94534.452
$ws.Range("N107").Value = -98374.452

# Row 131 (CUL)
$ws.Range("H131").Value = 3574523.2
$ws.Range("I131").Value = 930.5625
$ws.Range("K131").Value = 2791.6875
$ws.Range("M131").Value = 2248.3125

# Row 138 (CUL)
$ws.Range("H138").Value = 2486.5334
$ws.Range("I138").Value = 2413.077
$ws.Range("J138").Value = 2964
$ws.Range("K138").Value = 7239.231000000001
$ws.Range("L138").Value = 8892
$ws.Range("M138").Value = -2099.231000000001
$ws.Range("N138").Value = -19172

$ws = $wb.Worksheets.Item("GSM")
# Row 43 (GSM)
$ws.Range("H43").Value = 18287.777
$ws.Range("J43").Value = 23255.715
$ws.Range("L43").Value = 23255.715
$ws.Range("N43").Value = -23557.715

# Row 70 (GSM)
$ws.Range("H70").Value = 4798.4443
$ws.Range("I70").Value = 4775.4443
$ws.Range("J70").Value = 4844.4443
$ws.Range("K70").Value = 4775.4443
$ws.Range("L70").Value = 4844.4443
$ws.Range("M70").Value = -4505.4443
$ws.Range("N70").Value = -5384.4443

# Row 73 (GSM)
$ws.Range("H73").Value = 4798.4443
$ws.Range("I73").Value = 4775.4443
$ws.Range("J73").Value = 4844.4443
$ws.Range("K73").Value = 4775.4443
$ws.Range("L73").Value = 4844.4443
$ws.Range("M73").Value = -3839.4443
$ws.Range("N73").Value = -6716.4443

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 646.58826
$ws.Range("I46").Value = 671.4286
$ws.Range("J46").Value = 629.2
$ws.Range("K46").Value = 671.4286
$ws.Range("L46").Value = 629.2
$ws.Range("M46").Value = -483.4286
$ws.Range("N46").Value = -1005.2

# Row 47 (LTW)
$ws.Range("H47").Value = 14249
$ws.Range("J47").Value = 14249
$ws.Range("L47").Value = 14249
$ws.Range("N47").Value = -15229

# Row 52 (LTW)
$ws.Range("H52").Value = 14249
$ws.Range("J52").Value = 14249
$ws.Range("L52").Value = 14249
$ws.Range("N52").Value = -14715

# Row 74 (LTW)
$ws.Range("H74").Value = 34000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77 (LTW)
$ws.Range("H77").Value = 34000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 82 (WVR)
$ws.Range("H82").Value = 5050.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 5050.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 5050.5
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -5816.5

# Row 85 (WVR)
$ws.Range("H85").Value = 5050.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 5050.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 5050.5
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -7702.5
